# Applies the "Add files via upload" edit: adds a new KNN=500 results block
# (columns F:I, rows 31-35) mirroring the existing A:D block, and updates
# the original A:D block (rows 34-35) with new benchmark numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New header row 31/32 for the F:I block ---
# (apply the centered style before writing values so the engine reuses the
# existing "center" cell format instead of minting a duplicate one)
$ws.Range("F31:I31").HorizontalAlignment = -4108  # xlCenter
$ws.Range("F31").Value = "Avg MFCC (12 coeff) + Avg Delta + 500 knn"
$ws.Range("F31:I31").Merge() | Out-Null

$ws.Range("G32:I32").HorizontalAlignment = -4108  # xlCenter
$ws.Range("G32").Value = "EER"
$ws.Range("G32:I32").Merge() | Out-Null

# --- New sub-header row 33 for the F:I block ---
$ws.Range("F33").Value = "Train"
$ws.Range("G33").Value = "Test: Read"
$ws.Range("H33").Value = "Test: Phone"
$ws.Range("I33").Value = "Test: Mismatch"

# --- New data rows 34/35 for the F:I block (old A:D numbers moved here) ---
$ws.Range("F34").Value = "Read"
$ws.Range("G34").Value = 27.6423
$ws.Range("H34").Value = 28.3333
$ws.Range("I34").Value = 41.4815

$ws.Range("F35").Value = "Phone"
$ws.Range("G35").Value = 31.1111
$ws.Range("H35").Value = 21.4035
$ws.Range("I35").Value = 45.1852

# --- Updated data in the original A:D block ---
$ws.Range("B34").Value = 17.3984
$ws.Range("C34").Value = 28.3333
$ws.Range("D34").Value = 35.5556

$ws.Range("B35").Value = 33.3333
$ws.Range("C35").Value = 17.076
$ws.Range("D35").Value = 48.1481

# --- Selection shown when the file is saved ---
$ws.Range("E33").Select() | Out-Null
